$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (C and D) for "modelo" and "politica"
$ws.Columns("C:D").Insert()

# Header row
$ws.Cells.Item(1,3).Value = "modelo"
$ws.Cells.Item(1,4).Value = "politica"

# Row 2
$ws.Cells.Item(2,3).Value = "Sem Modelo"
$ws.Cells.Item(2,4).Value = ""
$ws.Cells.Item(2,6).Value = "premium"
$ws.Cells.Item(2,7).Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27685629?pdp_filters=seller_id:230239613#searchVariation=MLB27685629&position=3&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 3
$ws.Cells.Item(3,3).Value = "FONTE 120A"
$ws.Cells.Item(3,4).Value = "Igual"
$ws.Cells.Item(3,6).Value = "premium"
$ws.Cells.Item(3,7).Value = "https://www.mercadolivre.com.br/fonte-automotiva-120a-amperes-jfa-carregador-cor-preto/p/MLB21392652?pdp_filters=seller_id:230239613#searchVariation=MLB21392652&position=13&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 4
$ws.Cells.Item(4,3).Value = "FONTE 200 BOB"
$ws.Cells.Item(4,4).Value = "Igual"
$ws.Cells.Item(4,6).Value = "premium"
$ws.Cells.Item(4,7).Value = "https://www.mercadolivre.com.br/fonte-automotiva-jfa-storm-200a-bob-carregador-automatico-bivolt-cor-bob-200a-jfa/p/MLB24834408?pdp_filters=seller_id:230239613#searchVariation=MLB24834408&position=17&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 5
$ws.Cells.Item(5,3).Value = "FONTE 200 MONO"
$ws.Cells.Item(5,4).Value = "Acima"
$ws.Cells.Item(5,6).Value = "premium"
$ws.Cells.Item(5,7).Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-storm-voltimetro-digital-mono-220v-cor-preto/p/MLB24006449?pdp_filters=seller_id:230239613#searchVariation=MLB24006449&position=25&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 6
$ws.Cells.Item(6,3).Value = "Sem Modelo"
$ws.Cells.Item(6,4).Value = ""
$ws.Cells.Item(6,6).Value = "premium"
$ws.Cells.Item(6,7).Value = "https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-vermelho/p/MLB34210379?pdp_filters=seller_id:230239613#searchVariation=MLB34210379&position=30&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 7
$ws.Cells.Item(7,3).Value = "FONTE 120 BOB"
$ws.Cells.Item(7,4).Value = "Igual"
$ws.Cells.Item(7,6).Value = "premium"
$ws.Cells.Item(7,7).Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-120a-bob-slim-bivolt-cor-preto/p/MLB22144397?pdp_filters=seller_id:230239613#searchVariation=MLB22144397&position=32&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 8
$ws.Cells.Item(8,3).Value = "Modelo identificado mas fora do range de preco"
$ws.Cells.Item(8,4).Value = ""
$ws.Cells.Item(8,6).Value = "premium"
$ws.Cells.Item(8,7).Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-lite-storm-slim-bivolt-cor-azul/p/MLB24154371?pdp_filters=seller_id:230239613#searchVariation=MLB24154371&position=41&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 9
$ws.Cells.Item(9,3).Value = "Sem Modelo"
$ws.Cells.Item(9,4).Value = ""
$ws.Cells.Item(9,6).Value = "premium"
$ws.Cells.Item(9,7).Value = "https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-preto/p/MLB28687615?pdp_filters=seller_id:230239613#searchVariation=MLB28687615&position=45&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 10
$ws.Cells.Item(10,3).Value = "FONTE 90 BOB"
$ws.Cells.Item(10,4).Value = "Igual"
$ws.Cells.Item(10,6).Value = "premium"
$ws.Cells.Item(10,7).Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:230239613#searchVariation=MLB21562641&position=50&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 11
$ws.Cells.Item(11,3).Value = "Sem Modelo"
$ws.Cells.Item(11,4).Value = ""
$ws.Cells.Item(11,6).Value = "classico"
$ws.Cells.Item(11,7).Value = "https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-azul/p/MLB28722231?pdp_filters=seller_id:230239613#searchVariation=MLB28722231&position=31&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 12
$ws.Cells.Item(12,3).Value = "Modelo identificado mas fora do range de preco"
$ws.Cells.Item(12,4).Value = ""
$ws.Cells.Item(12,6).Value = "premium"
$ws.Cells.Item(12,7).Value = "https://www.mercadolivre.com.br/jfa-fonte-carregador-storm-lite-60a-3000-w-preto/p/MLB23456525?pdp_filters=seller_id:230239613#searchVariation=MLB23456525&position=33&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 13
$ws.Cells.Item(13,3).Value = "FONTE 70A"
$ws.Cells.Item(13,4).Value = "Igual"
$ws.Cells.Item(13,6).Value = "premium"
$ws.Cells.Item(13,7).Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-70a-bivolt-com-medidor-cca/p/MLB21455208?pdp_filters=seller_id:230239613#searchVariation=MLB21455208&position=35&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 14
$ws.Cells.Item(14,3).Value = "FONTE 60A"
$ws.Cells.Item(14,4).Value = "Acima"
$ws.Cells.Item(14,6).Value = "premium"
$ws.Cells.Item(14,7).Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca/p/MLB21320712?pdp_filters=seller_id:230239613#searchVariation=MLB21320712&position=38&search_layout=grid&type=product&tracking_id=eabc56ea-19f4-4216-a7bf-c1b6f3897ffc"

# Row 15
$ws.Cells.Item(15,3).Value = "Modelo identificado mas fora do range de preco"
$ws.Cells.Item(15,4).Value = ""
$ws.Cells.Item(15,6).Value = "premium"
$ws.Cells.Item(15,7).Value = "https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-120a-storm-lite-12v-bivolt-cor-preto/p/MLB23998473?pdp_filters=seller_id:230239613#searchVariation=MLB23998473&position=5&search_layout=grid&type=product&tracking_id=b4b39898-c93e-405a-af39-5c3e1d3eb3af"

# Row 16
$ws.Cells.Item(16,3).Value = "FONTE 40A"
$ws.Cells.Item(16,4).Value = "Acima"
$ws.Cells.Item(16,6).Value = "premium"
$ws.Cells.Item(16,7).Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-storm-40a-bivolt-12v-cor-preto/p/MLB22569833?pdp_filters=seller_id:230239613#searchVariation=MLB22569833&position=3&search_layout=grid&type=product&tracking_id=a17f50b6-be2d-4960-9c4c-53c711592ee9"

# Row 17
$ws.Cells.Item(17,3).Value = "Sem Modelo"
$ws.Cells.Item(17,4).Value = ""
$ws.Cells.Item(17,6).Value = "premium"
$ws.Cells.Item(17,7).Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:230239613#searchVariation=MLB27687422&position=42&search_layout=grid&type=product&tracking_id=a757e012-e79b-4160-bad3-ca3998c38d8c"

# Row 18
$ws.Cells.Item(18,3).Value = "Modelo identificado mas fora do range de preco"
$ws.Cells.Item(18,4).Value = ""
$ws.Cells.Item(18,6).Value = "premium"
$ws.Cells.Item(18,7).Value = "https://produto.mercadolivre.com.br/MLB-1948501526-carregador-bateria-portatil-carromoto-12v-jfa-f60a-redline-_JM#position%3D44%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D560c7adf-707a-4a2e-9939-b436d34088b1"

# Row 19
$ws.Cells.Item(19,3).Value = "Sem Modelo"
$ws.Cells.Item(19,4).Value = ""
$ws.Cells.Item(19,6).Value = "premium"
$ws.Cells.Item(19,7).Value = "https://produto.mercadolivre.com.br/MLB-1223504059-controle-longa-distancia-jfa-k1200-1200-mts-cor-azul-_JM#position%3D1%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3De7a8c97a-6d5e-46db-80ed-2e4af88ba9b4"

# Row 20
$ws.Cells.Item(20,3).Value = "Sem Modelo"
$ws.Cells.Item(20,4).Value = ""
$ws.Cells.Item(20,6).Value = "classico"
$ws.Cells.Item(20,7).Value = "https://produto.mercadolivre.com.br/MLB-1223662205-controle-longa-distancia-jfa-k1200-1200-mts-cor-vermelho-_JM#position%3D39%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D6efb3741-47a7-4fb4-8265-4554c5938651"

# Row 21
$ws.Cells.Item(21,3).Value = "Sem Modelo"
$ws.Cells.Item(21,4).Value = ""
$ws.Cells.Item(21,6).Value = "classico"
$ws.Cells.Item(21,7).Value = "https://produto.mercadolivre.com.br/MLB-1223499443-controle-longa-distancia-jfa-k1200-1200-mts-cor-preto-_JM#position%3D47%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D6efb3741-47a7-4fb4-8265-4554c5938651"
